$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -8.106000000000002
$ws.Range("A4").Value = -21.703
$ws.Range("B4").Value = 5.412000000000001
$ws.Range("A6").Value = -22.209
$ws.Range("A7").Value = -20.451
$ws.Range("C7").Value = -12.154
$ws.Range("C8").Value = -11.861
$ws.Range("B9").Value = 5.519
$ws.Range("C10").Value = -13.173
$ws.Range("B12").Value = 4.912999999999999
$ws.Range("C13").Value = -12.379
$ws.Range("D13").Value = -7.787000000000001
$ws.Range("A16").Value = -21.507
$ws.Range("C16").Value = -12.696
$ws.Range("B17").Value = 5.673999999999999
$ws.Range("B18").Value = 5.743
$ws.Range("B19").Value = 6.792
$ws.Range("A20").Value = -21.711
$ws.Range("B20").Value = 6.343999999999999
$ws.Range("D20").Value = -7.657000000000001
$ws.Range("D25").Value = -7.987
$ws.Range("B26").Value = 5.588
$ws.Range("A28").Value = -21.619
$ws.Range("A29").Value = -21.128
$ws.Range("C30").Value = -11.774
$ws.Range("B31").Value = 5.958
$ws.Range("A32").Value = -21.452
$ws.Range("D34").Value = -7.715999999999999
$ws.Range("B39").Value = 6.648999999999999
$ws.Range("D39").Value = -7.742
$ws.Range("A40").Value = -21.587
$ws.Range("B40").Value = 5.786
$ws.Range("C40").Value = -11.257
$ws.Range("B41").Value = 6.458000000000001
$ws.Range("B42").Value = 6.163
$ws.Range("B43").Value = 5.779999999999999
$ws.Range("C44").Value = -12.053
$ws.Range("A46").Value = -21.589
$ws.Range("B47").Value = 5.455
$ws.Range("B48").Value = 5.44
$ws.Range("A51").Value = -21.545
$ws.Range("D51").Value = -8.346
$ws.Range("A52").Value = -21.651
$ws.Range("A57").Value = -21.806
$ws.Range("A59").Value = -21.781
$ws.Range("D59").Value = -8.010000000000002
$ws.Range("D61").Value = -7.923
$ws.Range("A62").Value = -21.779
$ws.Range("B63").Value = 5.252
$ws.Range("B64").Value = 5.802
$ws.Range("D64").Value = -7.906999999999999
$ws.Range("A66").Value = -21.23
$ws.Range("A73").Value = -21.021
$ws.Range("A74").Value = -20.504
$ws.Range("B76").Value = 5.976000000000001
$ws.Range("D78").Value = -8.638
$ws.Range("B81").Value = 5.979
$ws.Range("D83").Value = -8.236999999999998
$ws.Range("B89").Value = 5.315
$ws.Range("C89").Value = -14.055
$ws.Range("C91").Value = -12.801
$ws.Range("A92").Value = -21.528
$ws.Range("D92").Value = -7.117
$ws.Range("B94").Value = 5.279999999999999
$ws.Range("D98").Value = -7.468999999999999
$ws.Range("A100").Value = -21.26
$ws.Range("D100").Value = -8.099
